# Update "想去人数" (interested-count) values on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages build 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3317   # was 3310
$wsExhibit.Range("F5").Value = 1333   # was 1328
$wsExhibit.Range("F6").Value = 316    # was 315

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3317   # was 3310
$wsAll.Range("F5").Value = 1333   # was 1328
$wsAll.Range("F7").Value = 316    # was 315
